$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Never went to school / never completed primary school (ISCED 0).global.lvl_educ"
$ws.Range("C1").Value = "Primary level of education (ISCED 1).global.lvl_educ"
$ws.Range("D1").Value = "Secondary level of education (ISCED 2 and ISCED 3).global.lvl_educ"
$ws.Range("E1").Value = "Higher education (ISCED 4 to ISCED 6).global.lvl_educ"
$ws.Range("F1").Value = "Not known / missing.global.lvl_educ"
$ws.Range("G1").Value = "Total.global.lvl_educ"
